$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.094.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5146"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.74%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3758"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07157"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8892"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.68"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.36%  "

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.901.53"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.02%  "

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07604"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.61%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.295"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.60"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.08%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9994"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.25%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008469"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.79%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.20%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.127.56"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.79%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.028"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.51%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.101.59"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.97%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.29%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.458"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.839"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.62%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.55"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.97"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.109"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.00%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.69"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.80%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.662"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.18%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.699"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.75%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09133"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05122"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.066"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.58%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.157"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.23%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7266"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.07%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02039"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.74%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.510"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.051"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.94%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.073"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5340"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.563"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.67"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.25%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.303"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.22%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1467"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.27%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4637"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.66%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9990"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.17%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.987"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.66%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.571"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.19%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.53"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.44%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.82"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.78%  "
